$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("TestSheet3")

# --- Rebuild TestSheet3 header row in the exact order needed so new shared
# strings get appended to sharedStrings.xml in the order the target
# workbook expects (Username, Id, Unemployment, Blank, Empty, then "  ") ---
$ws3.Range("B1").Value = "Username"
$ws3.Range("A1").Value = "Id"
$ws3.Range("C1").Value = "Date of Birth"
$ws3.Range("D1").Value = "Unemployment"
$ws3.Range("E1").Value = "Blank"
$ws3.Range("F1").Value = "Empty"

# Header row keeps the same bold "Heading 2" look as before, now across A:F.
$ws3.Range("A1:F1").Style = "Heading 2"

# --- Data rows ---
$ws3.Range("A2").Value = 1
$ws3.Range("B2").Value = "User5"
$ws3.Range("C2").Value = 30085
$ws3.Range("C2").NumberFormat = "mm-dd-yy"
$ws3.Range("D2").Value = $true

$ws3.Range("A3").Value = 2
$ws3.Range("B3").Value = "User6"
$ws3.Range("C3").Value = 28593
# Reuse C2's number-format style for C3 instead of minting a new one.
$ws3.Range("C2").Copy()
$ws3.Range("C3").PasteSpecial(-4122)
$ws3.Range("D3").Value = $false
$ws3.Range("E3").Value = "  "

# --- Column widths for the new layout ---
$ws3.Columns.Item(2).ColumnWidth = 19.5
$ws3.Columns.Item(3).ColumnWidth = 17.5
$ws3.Columns.Item(4).ColumnWidth = 16.5

# --- Page setup: portrait orientation ---
$ws3.PageSetup.Orientation = 1

# --- Selections ---
$ws1 = $wb.Worksheets.Item("TestSheet1")
$ws1.Range("A1:B3").Select()

$ws2 = $wb.Worksheets.Item("TestSheet2")
$ws2.Range("A1:B3").Select()

# TestSheet3 becomes the active sheet/tab with C2 selected
$ws3.Activate()
$ws3.Range("C2").Select()
